# Uren Registratie 12-04-2016 & Asset List Update & Menu
#
# Fills in the "Vrijdag" (Friday) hours for week 14 (row 78, Fahrettin's
# B-column total) and the whole week 15 block (rows 82-86), turning the
# previously-literal-zero weekly total (B87) back into a live SUM-style
# formula, matching the pattern used by every other week total in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 14 (rows 73-79): Vrijdag row picks up Fahrettin's 4 hours -------
$ws.Range("B78").Value = 4

# --- Week 15 (rows 81-87): fill in the whole week ---------------------
$ws.Range("B82").Value = 4

$ws.Range("B83").Value = 4
$ws.Range("C83").Value = 4
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = 4
$ws.Range("G83").Value = 4

$ws.Range("B84").Value = 2
$ws.Range("B85").Value = 2
$ws.Range("B86").Value = 4

# B87 was a stray literal 0 - restore it to the same kind of formula the
# other weekly-total cells (B71, B79, ...) use.
$ws.Range("B87").Formula = "=B82+B83+B84+B85+B86"

# --- View state: keep the on-screen selection in sync with the edits ------
$ws.Range("J77").Select()
